$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "% of Q Drop's" column (column I): header + one text value per
# course row. Written exactly like the other "% of ..." columns -- as
# literal text, e.g. "0.00%", "8.33%", "14.29%" -- not as numbers.
$qDropValues = @{
    "I1"  = "% of Q Drop's"
    "I3"  = "0.00%"
    "I6"  = "0.00%"
    "I9"  = "0.00%"
    "I12" = "0.00%"
    "I15" = "0.00%"
    "I18" = "0.00%"
    "I21" = "0.00%"
    "I24" = "0.00%"
    "I27" = "0.00%"
    "I30" = "8.33%"
    "I33" = "14.29%"
    "I36" = "0.00%"
}

foreach ($cellRef in $qDropValues.Keys) {
    $cell = $ws.Range($cellRef)
    # Briefly force text number-format so Excel stores the "NN.NN%"
    # strings verbatim instead of reinterpreting them as numeric
    # percentages, then drop the format again so the cell ends up with
    # no explicit style -- same as every other text cell on the sheet.
    $cell.NumberFormat = "@"
    $cell.Value = $qDropValues[$cellRef]
    $cell.ClearFormats()
}
